# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this week's play-by-play yardage logs to the running
# per-play sequences (Rush/Pass for OFF and DEF).
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value2 = $ydsWs.Range("B2").Value2 + " " + "6 0 2 4 5 5 6 5 1 -1 5 8 4 3 5 22 2 12 1 8 2 4 1 1 3"
$ydsWs.Range("B3").Value2 = $ydsWs.Range("B3").Value2 + " " + "6 2 6 17 2 2 9 8 4 19 3 56 4 4 15"
$ydsWs.Range("C2").Value2 = $ydsWs.Range("C2").Value2 + " " + "4 2 4 4 3 5 -1 5 2 2 1 3 2 6 7 2 2 5 6 14 4 5 7 1 0 -1 7 7 13 1 -2 10 1 2"
$ydsWs.Range("C3").Value2 = $ydsWs.Range("C3").Value2 + " " + "5 8 7 28 -1 10 3 11 -1 8 9 11 8 25 2 24 24 5"

# ---------------------------------------------------------------------------
# OFF sheet: updated down/distance + play totals (Home row 2 / Road row 3)
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value2 = 369
$offWs.Range("F2").Value2 = 102
$offWs.Range("G2").Value2 = 108
$offWs.Range("H2").Value2 = 5
$offWs.Range("I2").Value2 = 10
$offWs.Range("J2").Value2 = 52
$offWs.Range("N2").Value2 = 50

$offWs.Range("C3").Value2 = 315
$offWs.Range("E3").Value2 = 65
$offWs.Range("F3").Value2 = 189
$offWs.Range("G3").Value2 = 68
$offWs.Range("H3").Value2 = 62
$offWs.Range("I3").Value2 = 101
$offWs.Range("J3").Value2 = 108
$offWs.Range("L3").Value2 = 514
$offWs.Range("M3").Value2 = 331
$offWs.Range("Q3").Value2 = 925

# ---------------------------------------------------------------------------
# DEF sheet: updated down/distance + play totals (Home row 2 / Road row 3)
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value2 = 365
$defWs.Range("F2").Value2 = 133
$defWs.Range("G2").Value2 = 86
$defWs.Range("J2").Value2 = 66
$defWs.Range("N2").Value2 = 29
$defWs.Range("O2").Value2 = 32

$defWs.Range("B3").Value2 = 16
$defWs.Range("C3").Value2 = 355
$defWs.Range("E3").Value2 = 62
$defWs.Range("F3").Value2 = 204
$defWs.Range("G3").Value2 = 61
$defWs.Range("H3").Value2 = 54
$defWs.Range("I3").Value2 = 120
$defWs.Range("J3").Value2 = 102
$defWs.Range("L3").Value2 = 576
$defWs.Range("M3").Value2 = 368
$defWs.Range("Q3").Value2 = 982

# ---------------------------------------------------------------------------
# ST sheet: special-teams counts + appended kickoff/punt/return yardage logs
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value2 = 151
$stWs.Range("D2").Value2 = 126
$stWs.Range("H2").Value2 = 4
$stWs.Range("J2").Value2 = 49
$stWs.Range("K2").Value2 = 47
$stWs.Range("N2").Value2 = 20
$stWs.Range("O2").Value2 = 16
$stWs.Range("B3").Value2 = 106

$stWs.Range("B6").Value2 = $stWs.Range("B6").Value2 + " " + "23"
$stWs.Range("D3").Value2 = $stWs.Range("D3").Value2 + " " + "57 57 41 61 45 56 51"
$stWs.Range("D4").Value2 = $stWs.Range("D4").Value2 + " " + "0 0 0 0 0 12 18"
$stWs.Range("D5").Value2 = $stWs.Range("D5").Value2 + " " + "12 0 0 5 0 5"

# ---------------------------------------------------------------------------
# TURNS sheet: Road turnover counts
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("D3").Value2 = 16
$turnsWs.Range("E3").Value2 = 13

# ---------------------------------------------------------------------------
# PEN sheet: penalty counts
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value2 = 33
$penWs.Range("D2").Value2 = 13
$penWs.Range("B3").Value2 = 29
$penWs.Range("D4").Value2 = 16
